# Rebuild the data rows (A2:C17) of Sheet1 to reflect the refreshed
# false-positive feature summary (rdap_* features added, counts updated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Feature name, "Most Common SHAP Value", Count - in final row order (rows 2-17)
$rows = @(
    @("rdap_domain_age", 904408292000000000, 4),
    @("rdap_domain_active_time", 725760000000000000, 7),
    @("rdap_time_from_last_change", 400056294000000000, 2),
    @("tls_root_cert_lifetime", 5078, 82),
    @("tls_root_cert_validity_len", 1349, 7),
    @("rdap_ip_longest_v4_prefix_len", 25, 1),
    @("lex_name_len", 12, 1),
    @("ip_v4_count", 8, 21),
    @("lex_tld_len", 6, 2),
    @("dns_dn_level", 2, 12),
    @("ip_as_address_entropy", 2, 1),
    @("lex_sub_digit_ratio", 0.67, 1),
    @("dns_zone_entropy", 0.53, 6),
    @("dns_ttl_low", 0, 12),
    @("rdap_registrant_name_hash", 0, 7),
    @("rdap_registrar_name_hash", 0, 3)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
